$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New duty-roster values for column B (row -> name). $null clears the cell.
$values = @{
    2  = '遠藤隼人'
    3  = '富澤天音'
    4  = '神山修造'
    5  = '志塚惇希'
    6  = '川田涼介'
    7  = '豊島亮'
    8  = '兒島大志郎'
    9  = '日高泰聖'
    10 = '白岩詩佑介'
    11 = 'Cox Matthew Jonah'
    12 = 'Hansen Jakob U'
    13 = '石井海成'
    14 = 'Nicholas Tristan Aryasatyo'
    15 = $null
    16 = '小溝賢'
    17 = '小野文哉'
    18 = '渡部魁'
    19 = '崎谷航平'
    20 = '三神佳誠'
    21 = '氏家琉貴'
    22 = '羽賀尚生'
    23 = '島田実'
    24 = '足立耕平'
    25 = $null
    26 = '志塚惇希'
    27 = $null
    28 = '豊島亮'
    29 = $null
    30 = '日高泰聖'
    31 = $null
    32 = 'Cox Matthew Jonah'
}

foreach ($row in 2..32) {
    $val = $values[$row]
    if ($null -eq $val) {
        $ws.Range("B$row").Value = $null
    } else {
        $ws.Range("B$row").Value = $val
    }
}

$ws.Range("B13").Select()
